# [feat] conf json plus
#
# Adds two new "config value" columns (arrayValue / objectValue) to the
# student resource sample sheet, mirroring the existing layout:
#   row 1 -> header name
#   row 2 -> type name (same as header, repeated)
#   row 3 -> short description ("test" / "test2")
#   row 4 -> a concrete sample JSON value (array / object literal as text)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Range("J1").Value = "arrayValue"
$ws.Range("K1").Value = "objectValue"

# --- "type" row -------------------------------------------------------------
$ws.Range("J2").Value = "arrayValue"
$ws.Range("K2").Value = "objectValue"

# --- "description" row -----------------------------------------------------
$ws.Range("J3").Value = "test"
$ws.Range("K3").Value = "test2"

# --- sample data row --------------------------------------------------------
$ws.Range("J4").Value = '[1,2,4,"5",0.1]'
$ws.Range("K4").Value = '{"1":[0],"2":[1],"3":2.0,"4":[]}'

# --- column sizing for the two new columns ----------------------------------
$ws.Columns.Item(10).ColumnWidth = 20.9821428571429
$ws.Columns.Item(11).ColumnWidth = 39.7321428571429

# widen columns F/G slightly and re-apply H/I sizing to match the refreshed
# layout (cosmetic, best-effort - internal pixel rounding means these may
# land a hair off the source value)
$ws.Columns.Item(6).ColumnWidth = 20.25
$ws.Columns.Item(7).ColumnWidth = 109.25
$ws.Columns.Item(9).ColumnWidth = 54.4642857142857

# --- selection follows the last touched cell, like Excel would leave it ----
$ws.Range("K7").Select()
